# Update automatico via Actualizar 02-06-2021 19-09-47
#
# The DOLAR_OBS_ADO sheet is fed by a Banco Central de Chile web query
# (xl/connections.xml). This commit is the result of refreshing that query:
# 79 new daily rows are appended (serial dates 44147..44225, i.e.
# 2020-11-12 .. 2021-01-29 => sheet rows 684..762), the workbook-level
# "DOLAR_OBS_ADO" named range and the (hidden) AutoFilter database name are
# stretched to the new extent, and the view settles on the new last cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DOLAR_OBS_ADO")

# (date serial, value) for new rows 684..762. A number is the published
# "Dolar observado" rate; "--" marks a day the Central Bank did not publish
# one (weekends/holidays), matching the existing sharedString "--" (si #4).
$newRows = @(
    (44147,757.42),
    (44148,757.43),
    (44149,"--"),
    (44150,"--"),
    (44151,766.7),
    (44152,767.86),
    (44153,767.05),
    (44154,758.1),
    (44155,758.62),
    (44156,"--"),
    (44157,"--"),
    (44158,761.55),
    (44159,765.96),
    (44160,772.83),
    (44161,771.68),
    (44162,766),
    (44163,"--"),
    (44164,"--"),
    (44165,766.69),
    (44166,767.29),
    (44167,760.16),
    (44168,755.34),
    (44169,752.03),
    (44170,"--"),
    (44171,"--"),
    (44172,747.61),
    (44173,"--"),
    (44174,744.82),
    (44175,739.45),
    (44176,738.17),
    (44177,"--"),
    (44178,"--"),
    (44179,733.55),
    (44180,731.58),
    (44181,734.23),
    (44182,735.09),
    (44183,723.44),
    (44184,"--"),
    (44185,"--"),
    (44186,723.85),
    (44187,730.7),
    (44188,728.96),
    (44189,716.25),
    (44190,"--"),
    (44191,"--"),
    (44192,"--"),
    (44193,710.26),
    (44194,710.64),
    (44195,711.24),
    (44196,"--"),
    (44197,"--"),
    (44198,"--"),
    (44199,"--"),
    (44200,710.95),
    (44201,702.93),
    (44202,702.29),
    (44203,696.18),
    (44204,709.99),
    (44205,"--"),
    (44206,"--"),
    (44207,713.28),
    (44208,718.89),
    (44209,725.24),
    (44210,739.72),
    (44211,735.35),
    (44212,"--"),
    (44213,"--"),
    (44214,735.06),
    (44215,736.11),
    (44216,733.73),
    (44217,730.38),
    (44218,715.56),
    (44219,"--"),
    (44220,"--"),
    (44221,724.26),
    (44222,731.92),
    (44223,731),
    (44224,736.88),
    (44225,741.4)
)

$startRow = 684
# Style "templates" already present on the sheet: row 683 col B is a plain
# observed-rate number (xf s="3"); row 680 col B is the "--" not-published
# text (xf s="4"). Column A's date style (xf s="2") is stable either way.
$numTemplateRow = 683
$dashTemplateRow = 680

$r = $startRow
foreach ($pair in $newRows) {
    $serial = $pair[0]
    $val = $pair[1]

    $ws.Cells.Item($numTemplateRow, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial(-4122)
    $ws.Cells.Item($r, 1).Value = $serial

    if ("$val" -eq "--") {
        $ws.Cells.Item($dashTemplateRow, 2).Copy()
        $ws.Cells.Item($r, 2).PasteSpecial(-4122)
        $ws.Cells.Item($r, 2).Value = "--"
    } else {
        $ws.Cells.Item($numTemplateRow, 2).Copy()
        $ws.Cells.Item($r, 2).PasteSpecial(-4122)
        $ws.Cells.Item($r, 2).Value = $val
    }

    $r = $r + 1
}

$lastRow = $r - 1

# Stretch the workbook-level named range that mirrors the query's result area.
$dolarName = $wb.Names.Item("DOLAR_OBS_ADO")
$dolarName.RefersTo = "=DOLAR_OBS_ADO!`$A`$1:`$B`$" + $lastRow

# Re-create the hidden AutoFilter-database name over the refreshed data body
# (header row 3 through the new last row).
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=DOLAR_OBS_ADO!`$A`$3:`$B`$" + $lastRow)
$filterName.Visible = $false

# Column widths settle slightly narrower/wider after the refresh's autofit.
$ws.Columns.Item(1).ColumnWidth = 12.71
$ws.Columns.Item(2).ColumnWidth = 14.43

# Move the frozen-pane anchor and the active selection to track the newly
# appended last row/cell, as Excel does right after a query refresh.
$ws.Activate()
$ws.Cells.Item($lastRow, 2).Select()
